$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of the existing header row (copy from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Save column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
